# Append: 2025-09-29 06:36 JST
# Update the "取得日時" (retrieved datetime) column A for all existing
# data rows on the "ランサーズ" sheet from the previous run timestamp
# "2025-09-29 06:27:26" to the new one "2025-09-29 06:36:36".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-09-29 06:27:26"
$newTimestamp = "2025-09-29 06:36:36"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $current = $cell.Value()
    if ($current -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
